$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 42 and 43 swap their identity (ASV_ID / Species_name / Common_name / Category),
# and the J column value (Station18) moves from row 42 to row 43.

# Row 42 becomes the "Homo sapiens / Human" entry
$ws.Range("A42").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B42").Value = "Homo sapiens"
$ws.Range("C42").Value = "Human"
$ws.Range("D42").Value = "Human"
$ws.Range("J42").ClearContents()

# Row 43 becomes the "unassigned" entry
$ws.Range("A43").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B43").Value = "unassigned"
$ws.Range("C43").Value = "unassigned"
$ws.Range("D43").Value = "unassigned"
$ws.Range("J43").Value = 0
